$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 1051
$ws.Cells.Item(41, 9).Value = 907.5
$ws.Cells.Item(41, 10).Value = 1146.6666
$ws.Cells.Item(41, 11).Value = 907.5
$ws.Cells.Item(41, 12).Value = 1146.6666
$ws.Cells.Item(41, 13).Value = -467.5
$ws.Cells.Item(41, 14).Value = -2026.6666

$ws.Cells.Item(51, 8).Value = 1920.4
$ws.Cells.Item(51, 9).Value = 2000
$ws.Cells.Item(51, 10).Value = 1900.5
$ws.Cells.Item(51, 11).Value = 2000
$ws.Cells.Item(51, 12).Value = 1900.5
$ws.Cells.Item(51, 13).Value = -1516
$ws.Cells.Item(51, 14).Value = -2868.5

$ws.Cells.Item(76, 8).Value = 2981271.8
$ws.Cells.Item(76, 9).Value = 3551146.2
$ws.Cells.Item(76, 10).Value = 5260.1113
$ws.Cells.Item(76, 11).Value = 3551146.2
$ws.Cells.Item(76, 12).Value = 5260.1113
$ws.Cells.Item(76, 13).Value = -3550831.2
$ws.Cells.Item(76, 14).Value = -5890.1113

$ws.Cells.Item(79, 8).Value = 2981271.8
$ws.Cells.Item(79, 9).Value = 3551146.2
$ws.Cells.Item(79, 10).Value = 5260.1113
$ws.Cells.Item(79, 11).Value = 3551146.2
$ws.Cells.Item(79, 12).Value = 5260.1113
$ws.Cells.Item(79, 13).Value = -3550054.2
$ws.Cells.Item(79, 14).Value = -7444.1113

$ws.Cells.Item(132, 8).Value = 8773337
$ws.Cells.Item(132, 9).Value = 1118.862
$ws.Cells.Item(132, 11).Value = 3356.586
$ws.Cells.Item(132, 13).Value = -826.5860000000002

$ws.Cells.Item(137, 8).Value = 1482.2046
$ws.Cells.Item(137, 9).Value = 1227.069
$ws.Cells.Item(137, 11).Value = 3681.207
$ws.Cells.Item(137, 13).Value = -1131.207

$ws.Cells.Item(138, 8).Value = 2199.1604
$ws.Cells.Item(138, 9).Value = 776.69385
$ws.Cells.Item(138, 10).Value = 4377.3125
$ws.Cells.Item(138, 11).Value = 2330.08155
$ws.Cells.Item(138, 12).Value = 13131.9375
$ws.Cells.Item(138, 13).Value = 2809.91845
$ws.Cells.Item(138, 14).Value = -23411.9375

$ws.Cells.Item(141, 8).Value = 1523.6757
$ws.Cells.Item(141, 9).Value = 1110.5161
$ws.Cells.Item(141, 10).Value = 3658.3333
$ws.Cells.Item(141, 11).Value = 3331.5483
$ws.Cells.Item(141, 12).Value = 10974.9999
$ws.Cells.Item(141, 13).Value = 1848.4517
$ws.Cells.Item(141, 14).Value = -21334.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5256.7793
$ws.Cells.Item(32, 9).Value = 4119.258
$ws.Cells.Item(32, 10).Value = 9010.6
$ws.Cells.Item(32, 11).Value = 4119.258
$ws.Cells.Item(32, 12).Value = 9010.6
$ws.Cells.Item(32, 13).Value = -3832.258
$ws.Cells.Item(32, 14).Value = -9584.6

$ws.Cells.Item(74, 8).Value = 13890120
$ws.Cells.Item(74, 9).Value = 1035.2858
$ws.Cells.Item(74, 10).Value = 62501916
$ws.Cells.Item(74, 11).Value = 1035.2858
$ws.Cells.Item(74, 12).Value = 62501916
$ws.Cells.Item(74, 13).Value = -161.2858000000001
$ws.Cells.Item(74, 14).Value = -62503664

$ws.Cells.Item(77, 8).Value = 13890120
$ws.Cells.Item(77, 9).Value = 1035.2858
$ws.Cells.Item(77, 10).Value = 62501916
$ws.Cells.Item(77, 11).Value = 5176.429
$ws.Cells.Item(77, 12).Value = 312509580
$ws.Cells.Item(77, 13).Value = -808.4290000000001
$ws.Cells.Item(77, 14).Value = -312518316

$ws.Cells.Item(132, 8).Value = 2222.623
$ws.Cells.Item(132, 9).Value = 1617.0488
$ws.Cells.Item(132, 10).Value = 3464.05
$ws.Cells.Item(132, 11).Value = 4851.1464
$ws.Cells.Item(132, 12).Value = 10392.15
$ws.Cells.Item(132, 13).Value = -2321.1464
$ws.Cells.Item(132, 14).Value = -15452.15

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 1411.0834
$ws.Cells.Item(94, 9).Value = 962.1177
$ws.Cells.Item(94, 11).Value = 962.1177
$ws.Cells.Item(94, 13).Value = -511.1177

$ws.Cells.Item(134, 8).Value = 2941.4558
$ws.Cells.Item(134, 9).Value = 3291.8164
$ws.Cells.Item(134, 11).Value = 9875.449200000001
$ws.Cells.Item(134, 13).Value = -7340.449200000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2605164
$ws.Cells.Item(58, 9).Value = 3623736.2
$ws.Cells.Item(58, 10).Value = 2146.389
$ws.Cells.Item(58, 11).Value = 3623736.2
$ws.Cells.Item(58, 12).Value = 2146.389
$ws.Cells.Item(58, 13).Value = -3623533.2
$ws.Cells.Item(58, 14).Value = -2552.389

$ws.Cells.Item(132, 8).Value = 3450038.5
$ws.Cells.Item(132, 9).Value = 4546872
$ws.Cells.Item(132, 10).Value = 2848.2144
$ws.Cells.Item(132, 11).Value = 13640616
$ws.Cells.Item(132, 12).Value = 8544.643199999999
$ws.Cells.Item(132, 13).Value = -13638086
$ws.Cells.Item(132, 14).Value = -13604.6432

$ws.Cells.Item(134, 8).Value = 6946500
$ws.Cells.Item(134, 9).Value = 11496821
$ws.Cells.Item(134, 10).Value = 1272.421
$ws.Cells.Item(134, 11).Value = 34490463
$ws.Cells.Item(134, 12).Value = 3817.263
$ws.Cells.Item(134, 13).Value = -34487928
$ws.Cells.Item(134, 14).Value = -8887.262999999999

$ws.Cells.Item(136, 8).Value = 2605164
$ws.Cells.Item(136, 9).Value = 3623736.2
$ws.Cells.Item(136, 10).Value = 2146.389
$ws.Cells.Item(136, 11).Value = 10871208.6
$ws.Cells.Item(136, 12).Value = 6439.167
$ws.Cells.Item(136, 13).Value = -10868658.6
$ws.Cells.Item(136, 14).Value = -11539.167

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 17555520
$ws.Cells.Item(117, 9).Value = 20325.8
$ws.Cells.Item(117, 10).Value = 23818088
$ws.Cells.Item(117, 11).Value = 60977.39999999999
$ws.Cells.Item(117, 12).Value = 71454264
$ws.Cells.Item(117, 13).Value = -57535.39999999999
$ws.Cells.Item(117, 14).Value = -71461148

$ws.Cells.Item(133, 8).Value = 50263.523
$ws.Cells.Item(133, 9).Value = 88088.414
$ws.Cells.Item(133, 11).Value = 264265.242
$ws.Cells.Item(133, 13).Value = -259205.242

$ws.Cells.Item(134, 8).Value = 8984.333
$ws.Cells.Item(134, 9).Value = 9284.929
$ws.Cells.Item(134, 10).Value = 8816
$ws.Cells.Item(134, 11).Value = 27854.787
$ws.Cells.Item(134, 12).Value = 26448
$ws.Cells.Item(134, 13).Value = -22784.787
$ws.Cells.Item(134, 14).Value = -36588

$ws.Cells.Item(136, 8).Value = 17076.715
$ws.Cells.Item(136, 9).Value = 33845.668
$ws.Cells.Item(136, 10).Value = 4500
$ws.Cells.Item(136, 11).Value = 101537.004
$ws.Cells.Item(136, 12).Value = 13500
$ws.Cells.Item(136, 13).Value = -96437.00399999999
$ws.Cells.Item(136, 14).Value = -23700

$ws.Cells.Item(138, 8).Value = 7209.52
$ws.Cells.Item(138, 9).Value = 8838.056
$ws.Cells.Item(138, 10).Value = 3021.8572
$ws.Cells.Item(138, 11).Value = 26514.168
$ws.Cells.Item(138, 12).Value = 9065.5716
$ws.Cells.Item(138, 13).Value = -21374.168
$ws.Cells.Item(138, 14).Value = -19345.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 744.63635
$ws.Cells.Item(97, 9).Value = 640
$ws.Cells.Item(97, 10).Value = 927.75
$ws.Cells.Item(97, 11).Value = 640
$ws.Cells.Item(97, 12).Value = 927.75
$ws.Cells.Item(97, 13).Value = -144
$ws.Cells.Item(97, 14).Value = -1919.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 4168129.5
$ws.Cells.Item(22, 9).Value = 17857450
$ws.Cells.Item(22, 10).Value = 1814.4348
$ws.Cells.Item(22, 11).Value = 17857450
$ws.Cells.Item(22, 12).Value = 1814.4348
$ws.Cells.Item(22, 13).Value = -17857155
$ws.Cells.Item(22, 14).Value = -2404.4348

$ws.Cells.Item(27, 8).Value = 4168129.5
$ws.Cells.Item(27, 9).Value = 17857450
$ws.Cells.Item(27, 10).Value = 1814.4348
$ws.Cells.Item(27, 11).Value = 17857450
$ws.Cells.Item(27, 12).Value = 1814.4348
$ws.Cells.Item(27, 13).Value = -17857343
$ws.Cells.Item(27, 14).Value = -2028.4348

$ws.Cells.Item(46, 8).Value = 15152098
$ws.Cells.Item(46, 9).Value = 55556030
$ws.Cells.Item(46, 10).Value = 623.0625
$ws.Cells.Item(46, 11).Value = 55556030
$ws.Cells.Item(46, 12).Value = 623.0625
$ws.Cells.Item(46, 13).Value = -55555842
$ws.Cells.Item(46, 14).Value = -999.0625

$ws.Cells.Item(82, 8).Value = 357122.97
$ws.Cells.Item(82, 9).Value = 556886.8
$ws.Cells.Item(82, 10).Value = 80526.84
$ws.Cells.Item(82, 11).Value = 556886.8
$ws.Cells.Item(82, 12).Value = 80526.84
$ws.Cells.Item(82, 13).Value = -556525.8
$ws.Cells.Item(82, 14).Value = -81248.84

$ws.Cells.Item(85, 8).Value = 357122.97
$ws.Cells.Item(85, 9).Value = 556886.8
$ws.Cells.Item(85, 10).Value = 80526.84
$ws.Cells.Item(85, 11).Value = 556886.8
$ws.Cells.Item(85, 12).Value = 80526.84
$ws.Cells.Item(85, 13).Value = -555638.8
$ws.Cells.Item(85, 14).Value = -83022.84

$ws.Cells.Item(132, 8).Value = 10066508
$ws.Cells.Item(132, 9).Value = 12124892
$ws.Cells.Item(132, 10).Value = 3299.7778
$ws.Cells.Item(132, 11).Value = 36374676
$ws.Cells.Item(132, 12).Value = 9899.3334
$ws.Cells.Item(132, 13).Value = -36372146
$ws.Cells.Item(132, 14).Value = -14959.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 1350
$ws.Cells.Item(126, 9).Value = 550
$ws.Cells.Item(126, 10).Value = 1550
$ws.Cells.Item(126, 11).Value = 1650
$ws.Cells.Item(126, 12).Value = 4650
$ws.Cells.Item(126, 13).Value = 820
$ws.Cells.Item(126, 14).Value = -9590

$ws.Cells.Item(132, 8).Value = 925.43335
$ws.Cells.Item(132, 9).Value = 559.1395
$ws.Cells.Item(132, 10).Value = 1851.9412
$ws.Cells.Item(132, 11).Value = 1677.4185
$ws.Cells.Item(132, 12).Value = 5555.8236
$ws.Cells.Item(132, 13).Value = 852.5815
$ws.Cells.Item(132, 14).Value = -10615.8236

$ws.Cells.Item(136, 8).Value = 12581166
$ws.Cells.Item(136, 9).Value = 3203.625
$ws.Cells.Item(136, 11).Value = 9610.875
$ws.Cells.Item(136, 13).Value = -7060.875
